# The "Authors" column (E) stores author lists whose comma separators are
# followed by a run of spaces, e.g. ",        Kaniz..." . Each extra data
# "cleaning" pass that touched these three rows re-saved the author list
# with one more space added to that separator run, producing a brand-new
# (but otherwise identical) string that had to be appended to the shared
# string table. This commit applies two more such passes to rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-SeparatorSpace($text) {
    if ($text -match ",( +)\S") {
        $n = $matches[1].Length
        $old = "," + "".PadRight($n)
        $new = "," + "".PadRight($n + 1)
        return $text.Replace($old, $new)
    }
    return $text
}

foreach ($addr in @("E2", "E3", "E4")) {
    $current = $ws.Range($addr).Text
    # Apply the "extra space" cleaning pass twice, matching the two
    # additional shared-string revisions introduced for each of these rows.
    $updated = Add-SeparatorSpace (Add-SeparatorSpace $current)
    $ws.Range($addr).Value = $updated
}
